# Add new Harris Interactive poll rows (10/13 wave) to the bottom of the data table,
# warning: overestimated sample size, to be corrected when true sample size is available.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the three new respondent rows (55-57), one line per candidate combination
# scenario (turnout/second-round hypothesis), matching the layout already used for
# earlier Harris polls further up the sheet.
$newRows = @(
    @{ Row = 55; Values = @{
        A = 16; B = 2021; C = 7; D = 10; E = 10; F = "harris"; G = "online"; H = "included";
        I = 1051; J = 1; K = 1; L = 11; M = 2; N = 2; O = 7; P = 5; Q = 24;
        T = 14; U = "T_0.5"; V = 1; W = 15; X = 17;
        Y = "T_0.5"; Z = "T_0.5"; AA = "T_0.5"; AB = "T_0.5"
    } },
    @{ Row = 56; Values = @{
        A = 16; B = 2021; C = 7; D = 10; E = 10; F = "harris"; G = "online"; H = "included";
        I = 1051; J = 1; K = 1; L = 11; M = 2; N = 2; O = 8; P = 5; Q = 25;
        R = 11;
        U = "T_0.5"; V = 2; W = 15; X = 17;
        Y = "T_0.5"; Z = "T_0.5"; AA = "T_0.5"; AB = "T_0.5"
    } },
    @{ Row = 57; Values = @{
        A = 16; B = 2021; C = 7; D = 10; E = 10; F = "harris"; G = "online"; H = "included";
        I = 1051; J = 1; K = 1; L = 11; M = 2; N = 2; O = 8; P = 5; Q = 27;
        S = 7;
        U = "T_0.5"; V = 2; W = 16; X = 18;
        Y = "T_0.5"; Z = "T_0.5"; AA = "T_0.5"; AB = "T_0.5"
    } }
)

# Columns that, in the existing sheet, are styled with a plain black font (style index 1)
# to flag them as placeholder/"T_0.5" scenario columns.
$blackFontCols = @("Y", "Z", "AA", "AB")

foreach ($rowDef in $newRows) {
    $r = $rowDef.Row
    foreach ($col in $rowDef.Values.Keys) {
        $cell = $ws.Range("$col$r")
        $cell.Value = $rowDef.Values[$col]
        if ($blackFontCols -contains $col) {
            $cell.Font.Color = 0
        }
    }
}

# Update the active selection to reflect the newly added rows, as in the saved
# workbook (move the active cell to the new last data row).
$null = $ws.Range("P55").Select()
